# Insert a new price record for "Mandarina / Murcott / Primera" at row 350
# of the daily logic subset sheet. This pushes the existing rows 350-405
# down to 351-406 (dimension grows from A1:T405 to A1:T406) and fills the
# newly opened row 350 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 350:405 down by one row, opening up a blank row 350.
$ws.Rows("350:350").Insert()

# Populate the new row 350 with the new record.
$ws.Range("A350").Value = 5
$ws.Range("B350").Value = "Macroferia Regional de Talca"
$ws.Range("C350").Value = "Maule"
$ws.Range("D350").Value = 44964
$ws.Range("E350").Value = 7
$ws.Range("F350").Value = "Fruta"
$ws.Range("G350").Value = 100102
$ws.Range("H350").Value = "Cítricos"
$ws.Range("I350").Value = 100102004
$ws.Range("J350").Value = "Mandarina"
$ws.Range("K350").Value = "Murcott"
$ws.Range("L350").Value = "Primera"
$ws.Range("M350").Value = 300
$ws.Range("N350").Value = 10000
$ws.Range("O350").Value = 10000
$ws.Range("P350").Value = 10000
$ws.Range("Q350").Value = "`$/caja 15 kilos granel"
$ws.Range("R350").Value = "Región de Coquimbo"
$ws.Range("S350").Value = 667
$ws.Range("T350").Value = 15
